# Apply the "SAN# logic" update to EUC_Perth_Assets.xlsx
#
# This mirrors the behaviour described in the commit message: a handful of
# new inventory-count log entries were recorded (some of them for laptop /
# mini-pc items, which now also capture a "SAN" asset number), the rolling
# item-count tables were refreshed with newer totals, and the Build Room
# ("BR_*") sheets picked up their own running counts as items started moving
# through that workflow too.

$wb = $excel.ActiveWorkbook

# Helper: force a cell that has no value to still be materialised in the
# saved worksheet (as an empty cell) without pulling in any non-default
# formatting. Re-asserting the sheet's own default font name is enough to
# "touch" the cell while staying at style index 0.
function Touch-EmptyCell($range) {
    $range.Font.Name = "Arial"
}

# ---------------------------------------------------------------------
# Sheet "4.2_Items" - refreshed LastCount / NewCount totals
# ---------------------------------------------------------------------
$items42 = $wb.Worksheets.Item("4.2_Items")

$items42.Range("B2").Value = 1153
$items42.Range("C2").Value = 1154

$items42.Range("B3").Value = 1409
$items42.Range("C3").Value = 1411

$items42.Range("B4").Value = 940
$items42.Range("C4").Value = 990

$items42.Range("B5").Value = 96
$items42.Range("C5").Value = 99

$items42.Range("C7").Value = 47

$items42.Range("B8").Value = 78
$items42.Range("C8").Value = 108

$items42.Range("B11").Value = 78040
$items42.Range("C11").Value = 78020

# ---------------------------------------------------------------------
# Sheet "4.2_Timestamps" - new action log rows
# ---------------------------------------------------------------------
$ts42 = $wb.Worksheets.Item("4.2_Timestamps")

# Rows 47 & 48 already existed but picked up the standard row style used by
# the rest of the log (same format as row 45/46 above them).
$logStyle = $ts42.Range("A45").Style
$ts42.Range("A47:D48").Style = $logStyle

$monitor24 = "Monitor 24" + [char]0x201D + [char]0x00A0

# Row 49: Monitor 24" stock correction, logged with the usual row style and
# an (empty) SAN Number cell, same shape as rows 45-48.
$ts42.Range("A49").Value = "2023-12-08 23:35:08"
$ts42.Range("B49").Value = $monitor24
$ts42.Range("C49").Value = "Add 5000"
$ts42.Range("A49:D49").Style = $logStyle

# Rows 50 & 51: further Monitor 24" corrections - no SAN Number column used.
$ts42.Range("A50").Value = "2023-12-08 23:45:46"
$ts42.Range("B50").Value = $monitor24
$ts42.Range("C50").Value = "Subtract 68"
$ts42.Range("A50:C50").Style = $logStyle

$ts42.Range("A51").Value = "2023-12-08 23:46:04"
$ts42.Range("B51").Value = $monitor24
$ts42.Range("C51").Value = "Subtract 1"
$ts42.Range("A51:C51").Style = $logStyle

# Rows 52-56: next day's entries, written with default (un-styled) cells,
# same as the very first log rows in the sheet. Dock/Desktop Mini entries
# are not laptops or mini-PCs in the SAN sense, so their SAN cell is left
# blank (but still present). Laptop entries prompt for - and record - a SAN
# Number.
$ts42.Range("A52").Value = "2023-12-09 21:43:21"
$ts42.Range("B52").Value = "Dock Thunderbolt G4"
$ts42.Range("C52").Value = "Add 2"
Touch-EmptyCell $ts42.Range("D52")

$ts42.Range("A53").Value = "2023-12-09 21:43:25"
$ts42.Range("B53").Value = "Desktop Mini"
$ts42.Range("C53").Value = "Add 2"
Touch-EmptyCell $ts42.Range("D53")

$ts42.Range("A54").Value = "2023-12-09 21:43:31"
$ts42.Range("B54").Value = "Desktop Mini"
$ts42.Range("C54").Value = "Add 1"
Touch-EmptyCell $ts42.Range("D54")

$ts42.Range("A55").Value = "2023-12-09 21:43:49"
$ts42.Range("B55").Value = "Laptop 840 G9"
$ts42.Range("C55").Value = "Add 1"
$ts42.Range("D55").Value = "SAN123456"

$ts42.Range("A56").Value = "2023-12-09 21:45:37"
$ts42.Range("B56").Value = "Laptop 840 G9"
$ts42.Range("C56").Value = "Add 3"
$ts42.Range("D56").Value = "SAN111"

# ---------------------------------------------------------------------
# Sheet "BR_Items" - Build Room running counts now tracked per item
# ---------------------------------------------------------------------
$brItems = $wb.Worksheets.Item("BR_Items")
$brItemStyle = $items42.Range("B2").Style

$brItems.Range("B2").Value = 60
$brItems.Range("C2").Value = 80
$brItems.Range("B2:C2").Style = $brItemStyle

$brItems.Range("B3").Value = 0
$brItems.Range("C3").Value = 20
$brItems.Range("B3:C3").Style = $brItemStyle

$brItems.Range("B4").Value = 0
$brItems.Range("C4").Value = 20
$brItems.Range("B4:C4").Style = $brItemStyle

$brItems.Range("B5").Value = 0
$brItems.Range("C5").Value = 20
$brItems.Range("B5:C5").Style = $brItemStyle

$brItems.Range("B6").Value = 0
$brItems.Range("C6").Value = 20
$brItems.Range("B6:C6").Style = $brItemStyle

$brItems.Range("B7").Value = 0
$brItems.Range("C7").Value = 20
$brItems.Range("B7:C7").Style = $brItemStyle

$brItems.Range("C8").Value = 11143
$brItems.Range("C8").Style = $brItemStyle

# ---------------------------------------------------------------------
# Sheet "BR_Timestamps" - new Build Room action log rows
# ---------------------------------------------------------------------
$tsBR = $wb.Worksheets.Item("BR_Timestamps")
$logStyleBR = $tsBR.Range("A2").Style

# Row 5 already existed but now picks up the standard row style.
$tsBR.Range("A5:C5").Style = $logStyleBR

$tsBR.Range("A6").Value = "2023-12-08 23:35:20"
$tsBR.Range("B6").Value = "Laptop 840 G9"
$tsBR.Range("C6").Value = "Subtract 5000"
$tsBR.Range("A6:C6").Style = $logStyleBR

$tsBR.Range("A7").Value = "2023-12-08 23:35:23"
$tsBR.Range("B7").Value = "Laptop 840 G10"
$tsBR.Range("C7").Value = "Add 5000"
$tsBR.Range("A7:D7").Style = $logStyleBR

$tsBR.Range("A8").Value = "2023-12-08 23:45:11"
$tsBR.Range("B8").Value = $monitor24
$tsBR.Range("C8").Value = "Subtract 23"
$tsBR.Range("A8:C8").Style = $logStyleBR

$tsBR.Range("A9").Value = "2023-12-08 23:45:29"
$tsBR.Range("B9").Value = $monitor24
$tsBR.Range("C9").Value = "Subtract 68"
$tsBR.Range("A9:C9").Style = $logStyleBR

Write-Output "SAN# logic update applied"
